$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F11 (84 -> 85) and F14 (190 -> 191)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 85
$ws1.Range("F14").Value = 191

# Sheet "全部类型" - mirrors the same data, apply identical updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 85
$ws4.Range("F14").Value = 191
